# Remove the "Ver no Jupiter..." footer block (an empty paragraph followed
# by the "Ver no Jupiter Salvar em pdf Salvar em docx" line and the
# "(c) 2020 ... Creative Commons Attribution" copyright line) that used to
# follow the "Requisitos" section's last line.

$d = $word.ActiveDocument

$anchorText = "LOQ4064: Engenharia de Processos Quimicos I (Requisito fraco)"

$anchorPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd("`r") -eq $anchorText) {
        $anchorPara = $p
        break
    }
}

if ($anchorPara -ne $null) {
    # The three paragraphs right after the anchor are, in order:
    #   1) a blank "Normal" paragraph
    #   2) "Ver no Jupiter Salvar em pdf Salvar em docx"
    #   3) "© 2020 . Contact: luizeleno@usp.br. ... Creative Commons Attribution"
    $p1 = $anchorPara.Next()
    $p2 = $p1.Next()
    $p3 = $p2.Next()

    $deleteRange = $d.Range($p1.Range.Start, $p3.Range.End)
    $deleteRange.Delete()
}
